$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Resize the table and its columns (dxa -> points, 20 dxa = 1 pt) ---
# Overall table width: 8580 -> 7104 dxa
$t.PreferredWidth = 7104 / 20
# Switch to a fixed table layout (matches <w:tblLayout w:type="fixed"/>)
$t.AllowAutoFit = $false

# Column widths (dxa): 2580,1500,1500,1500,1500 -> 2221,1194,1350,979,1360
$newWidths = @(2221, 1194, 1350, 979, 1360)
for ($i = 1; $i -le $t.Columns.Count; $i++) {
    $t.Columns.Item($i).Width = $newWidths[$i - 1] / 20
}

# --- Split "grandis" into "grandi" + "s", with a _GoBack bookmark in between ---
$hit = $d.Content
$found = $hit.Find.Execute("grandis")
if ($found) {
    $splitPoint = $hit.Start + 6
    $insertionPoint = $d.Range($splitPoint, $splitPoint)
    $d.Bookmarks.Add("_GoBack", $insertionPoint) | Out-Null
}
